$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.848.23"
$ws.Range("E2").Value = "  +0.71%  "

$ws.Range("D3").Value = "2.655.68"
$ws.Range("E3").Value = "  +2.05%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'538.86"
$ws.Range("E5").Value = "  +0.58%  "

$ws.Range("D6").Value = "'146.53"
$ws.Range("E6").Value = "  +3.96%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +1.20%  "

$ws.Range("D9").Value = "2.670.48"
$ws.Range("E9").Value = "  +2.17%  "

$ws.Range("D10").Value = "'6.66"
$ws.Range("E10").Value = "  +2.72%  "

$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = "  +0.63%  "

$ws.Range("E12").Value = "  +0.85%  "

$ws.Range("E13").Value = "  -0.59%  "

$ws.Range("D14").Value = "3.124.13"
$ws.Range("E14").Value = "  +2.06%  "

$ws.Range("D15").Value = "59.754.13"
$ws.Range("E15").Value = "  +0.68%  "

$ws.Range("D16").Value = "'21.29"
$ws.Range("E16").Value = "  +3.41%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000136"
$ws.Range("E17").Value = "  +1.43%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.633.61"
$ws.Range("E18").Value = "  +0.70%  "

$ws.Range("D19").Value = "'344.57"
$ws.Range("E19").Value = "  -0.25%  "

$ws.Range("E20").Value = "  +1.95%  "

$ws.Range("D21").Value = "'10.45"
$ws.Range("E21").Value = "  +3.19%  "

$ws.Range("D22").Value = "'6.35"
$ws.Range("E22").Value = "  -0.35%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").Value = "'66.80"
$ws.Range("E24").Value = "  -0.53%  "

$ws.Range("D25").Value = "'0.418"
$ws.Range("E25").Value = "  +2.36%  "

$ws.Range("E26").Value = "  -0.85%  "

$ws.Range("E27").Value = "  +0.04%  "

$ws.Range("D28").Value = "'7.35"
$ws.Range("E28").Value = "  +1.75%  "

$ws.Range("D29").Value = "0.0₃0759"
$ws.Range("E29").Value = "  +1.66%  "

$ws.Range("E30").Value = "  -0.05%  "

$ws.Range("E31").Value = "  +1.60%  "

$ws.Range("D32").Value = "'5.87"
$ws.Range("E32").Value = "  +0.38%  "

$ws.Range("D33").Value = "'19.02"
$ws.Range("E33").Value = "  +0.82%  "

$ws.Range("D34").Value = "'150.61"
$ws.Range("E34").Value = "  +0.89%  "

$ws.Range("D35").Value = "'4.05"
$ws.Range("E35").Value = "  +1.17%  "

$ws.Range("E36").Value = "  +2.19%  "

$ws.Range("D37").Value = "'0.846"

$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("D39").Value = "'0.846"
$ws.Range("E39").Value = "  +1.09%  "

$ws.Range("D40").Value = "'291.78"
$ws.Range("E40").Value = "  +5.43%  "

$ws.Range("D41").Value = "'3.60"
$ws.Range("E41").Value = "  +1.66%  "

$ws.Range("D43").Value = "'0.609"
$ws.Range("E43").Value = "  +1.40%  "

$ws.Range("D44").Value = "'19.51"
$ws.Range("E44").Value = "  +4.76%  "

$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0540"
$ws.Range("E45").Value = "  +3.05%  "

$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "'10.72"
$ws.Range("E46").Value = "  -0.17%  "

$ws.Range("D47").Value = "'0.0952"
$ws.Range("E47").Value = "  -1.02%  "

$ws.Range("D48").Value = "1.983.00"
$ws.Range("E48").Value = "  +1.67%  "

$ws.Range("D49").Value = "'0.0228"
$ws.Range("E49").Value = "  +1.93%  "

$ws.Range("D50").Value = "'4.60"
$ws.Range("E50").Value = "  -2.43%  "

$ws.Range("D51").Value = "'18.47"
$ws.Range("E51").Value = "  +0.26%  "
